# edit.ps1
# Applies the "Updated cryptos list" refresh: rewrites the Price / Volume(1h) figures
# for every existing coin row, and replaces USDe (which fell out of the top list) with
# InjectiveProtocol/Stellar/ThetaToken/Cronos sliding up one slot and THORChain newly
# appearing at the bottom (row 51). The numeric rank in column A is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Writes $value into $range while forcing Excel to keep it as literal text
    # (prevents things like "2.00" or "23.80" from being silently turned into the
    # numbers 2 / 23.8 and losing their original formatting), exactly like a user
    # typing an apostrophe before a numeric-looking entry in the Excel UI.
    $range.Value = "'" + $value
}

# --- Update Price (D) / Volume(1h) (E) for the coins whose rank and identity do not change ---
Set-TextValue $ws.Range("D2") "60.917.52"
$ws.Range("E2").Value = "  +0.29%  "
Set-TextValue $ws.Range("D3") "2.919.33"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +1.22%  "
Set-TextValue $ws.Range("D6") "146.45"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  +0.07%  "
Set-TextValue $ws.Range("D8") "0.508"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("E14").Value = "  +0.18%  "
Set-TextValue $ws.Range("D15") "3.401.44"
$ws.Range("E15").Value = "  +0.31%  "
Set-TextValue $ws.Range("D16") "60.828.11"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("E17").Value = "  -0.53%  "
Set-TextValue $ws.Range("D18") "2.917.17"
$ws.Range("E18").Value = "  +0.22%  "
Set-TextValue $ws.Range("D19") "432.44"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("E22").Value = "  -0.92%  "
Set-TextValue $ws.Range("D23") "81.37"
$ws.Range("E23").Value = "  +1.23%  "
Set-TextValue $ws.Range("D24") "10.95"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -2.57%  "
Set-TextValue $ws.Range("D31") "26.69"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("E33").Value = "  -0.06%  "
Set-TextValue $ws.Range("D34") "0.0₃0858"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("E41").Value = "  -4.38%  "
Set-TextValue $ws.Range("D42") "40.18"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("E43").Value = "  +1.61%  "
Set-TextValue $ws.Range("D44") "2.692.49"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("E45").Value = "  -1.42%  "
Set-TextValue $ws.Range("D46") "133.36"
$ws.Range("E46").Value = "  +1.03%  "

# --- USDe drops out of the table; InjectiveProtocol/Stellar/ThetaToken/Cronos each move up one
#     row (47-50) keeping their existing rank number in column A, and THORChain is newly added
#     as row 51 ---
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "23.80"
$ws.Range("E47").Value = "  -2.01%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D48") "0.106"
$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D49") "2.00"
$ws.Range("E49").Value = "  -3.46%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.124"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "5.90"
$ws.Range("E51").Value = "  -1.43%  "
